$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "AddCustomerTest"
$ws.Range("A1").Value = "firstname"
$ws.Range("B1").Value = "lastname"
$ws.Range("C1").Value = "postcode"
$ws.Range("A2").Value = "Sagrika"
$ws.Range("B2").Value = "Srivastava"
$ws.Range("C2").Value = "'001100"

$ws.Columns("F:F").Group()
$ws.Columns("F:F").Group()
$ws.Columns("F:F").Delete()
$ws.Rows("5:5").Group()
$ws.Rows("5:5").Delete()

$ws.Range("N9").Select()
